# Retraining the models for 3D Steel and Elnet
# Shift the timestamp column (A) forward by one day (rows 2-97) and
# replace the "Actual Production (MW)" values (B) for rows 2-42 with
# the newly retrained figures. Rows 43-97 keep their value of 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @{
    2 = 1312
    3 = 1333
    4 = 1349
    5 = 1383
    6 = 1401
    7 = 1428
    8 = 1421
    9 = 1438
    10 = 1400
    11 = 1355
    12 = 1333
    13 = 1341
    14 = 1312
    15 = 1273
    16 = 1245
    17 = 1240
    18 = 1203
    19 = 1198
    20 = 1190
    21 = 1173
    22 = 1108
    23 = 1036
    24 = 1017
    25 = 1009
    26 = 960
    27 = 971
    28 = 978
    29 = 931
    30 = 836
    31 = 794
    32 = 731
    33 = 641
    34 = 546
    35 = 486
    36 = 429
    37 = 425
    38 = 438
    39 = 452
    40 = 466
    41 = 486
    42 = 524
}

for ($row = 2; $row -le 97; $row++) {
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $aCell.Value2() + 1

    if ($newB.ContainsKey($row)) {
        $ws.Cells.Item($row, 2).Value = $newB[$row]
    }
}
